$wb = $excel.ActiveWorkbook

# Sheet "展览" (Sheet1): increment F column "want to go" counts by 1 for specific rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 517
$wsExhibit.Range("F6").Value = 365
$wsExhibit.Range("F7").Value = 226
$wsExhibit.Range("F9").Value = 377
$wsExhibit.Range("F10").Value = 5452
$wsExhibit.Range("F11").Value = 126

# Sheet "演出" (Sheet2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 52

# Sheet "全部类型" (Sheet4): same events, combined list
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 52
$wsAll.Range("F4").Value = 517
$wsAll.Range("F7").Value = 365
$wsAll.Range("F8").Value = 226
$wsAll.Range("F12").Value = 377
$wsAll.Range("F13").Value = 5452
$wsAll.Range("F14").Value = 126
